$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# --- Overview sheet: row 3 is the "b.md" file, now ready for handoff ---
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"
$wsOverview.Range("D3").Value = "2016-27-19 00:27:12"

# --- zh-cn sheet: row 3 ("b.md") gets a new handoff file + status ---
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("D3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$wsZhCn.Range("E3").Value = "2016-03-19 00:27:09"

foreach ($hl in $wsZhCn.Hyperlinks) {
    if ($hl.Range.Address() -eq '$D$3') {
        $hl.TextToDisplay = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
    }
}

# --- de-de sheet: row 3 ("b.md") gets a new handoff file + status ---
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("D3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$wsDeDe.Range("E3").Value = "2016-03-19 00:27:12"

foreach ($hl in $wsDeDe.Hyperlinks) {
    if ($hl.Range.Address() -eq '$D$3') {
        $hl.TextToDisplay = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
    }
}
